# Insert a new column before column A, shifting all existing data right
# by one column (A->B, B->C, C->D, D->E, E->F), then populate the new
# column A with a zero-based segment index (and header "segments" moves
# into B1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

$ws.Columns.Item(1).Insert()

# Header for the new column - match the style of the other header cells.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial($xlPasteFormats)
$ws.Range("B1").Value = "segments"

# Fill the new index column (A2:A20) with 0-based row index, matching the
# style used by the (now shifted) segment-name column (B, which carries
# the original bold/bordered/centered style).
$ws.Range("B2").Copy()
$ws.Range("A2:A20").PasteSpecial($xlPasteFormats)

for ($i = 2; $i -le 20; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 2
}

$excel.CutCopyMode = 0
